# fix: revisar form antes de crear credenciales
#
# Adds a "Genero" / "M" question-answer pair as a new row (33) in the
# survey data, duplicating the fixed (non question/answer) columns from
# the row above it, and moves the active selection to H35 (previously J35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 currently has no data (blank record). Populate its "fixed" columns
# (A:O -- area/region/.../categoria) by copying them from the row above
# (row 32), which carries the same survey/beneficiary record. Using
# PasteSpecial(values) instead of setting .Value cell-by-cell keeps the
# original cell formatting/style (s="3") and avoids Excel re-typing text
# such as "03" as a number.
$ws.Range("A32:O32").Copy()
$ws.Range("A33:O33").PasteSpecial(-4163)

# New question/answer pair for row 33: Genero = M
$ws.Cells.Item(33, 16).Value = "Genero"
$ws.Cells.Item(33, 17).Value = "M"

# Move the active cell/selection from J35 to H35.
$ws.Range("H35").Select()
